# Cleaning up the incompatible types mess
# - Drops two stray empty placeholder cells left over in row 1303 (F/G),
#   and appends the new pick-list rows 1304:1335 that were missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1303 previously ended with two empty leftover cells in F/G; clear them.
$ws.Range("F1303").ClearContents()
$ws.Range("G1303").ClearContents()

# New order rows pulled in from the latest invoices (Feb 17 - Feb 26, 2025).
# Row 1304
$ws.Range("A1304").Value2 = ' February 17, 2025'
$ws.Range("B1304").Value2 = '''107750625'
$ws.Range("C1304").Value2 = '110TWN09180'
$ws.Range("D1304").Value2 = 'Twinings Lemon & Ginger Herbal Tea Bags - 25/Box'
$ws.Range("E1304").Value2 = ''' 12'
$ws.Range("H1304").Value2 = 4.59

# Row 1305
$ws.Range("A1305").Value2 = ' February 17, 2025'
$ws.Range("B1305").Value2 = '''107750625'
$ws.Range("C1305").Value2 = '110TWN09181'
$ws.Range("D1305").Value2 = 'Twinings English Breakfast Tea Bags - 25/Box'
$ws.Range("E1305").Value2 = ''' 24'
$ws.Range("H1305").Value2 = 4.59

# Row 1306
$ws.Range("A1306").Value2 = ' February 17, 2025'
$ws.Range("B1306").Value2 = '''107750625'
$ws.Range("C1306").Value2 = '110TWN09183'
$ws.Range("D1306").Value2 = 'Twinings Earl Grey Tea Bags - 25/Box'
$ws.Range("E1306").Value2 = ''' 6'
$ws.Range("H1306").Value2 = 4.59

# Row 1307
$ws.Range("A1307").Value2 = ' February 17, 2025'
$ws.Range("B1307").Value2 = '''107750625'
$ws.Range("C1307").Value2 = '''40862028'
$ws.Range("D1307").Value2 = 'Ghirardelli 30 lb. Sweet Ground Chocolate & Cocoa Powder'
$ws.Range("E1307").Value2 = ''' 5'
$ws.Range("H1307").Value2 = 123.47

# Row 1308
$ws.Range("A1308").Value2 = ' February 17, 2025'
$ws.Range("B1308").Value2 = '''107750625'
$ws.Range("C1308").Value2 = '8808604CS'
$ws.Range("D1308").Value2 = 'Torani Puremade White Chocolate Flavoring Sauce 64 fl. oz. - 4/Case'
$ws.Range("E1308").Value2 = ''' 1'
$ws.Range("H1308").Value2 = 72.98999999999999

# Row 1309
$ws.Range("A1309").Value2 = ' February 17, 2025'
$ws.Range("B1309").Value2 = '''107750625'
$ws.Range("C1309").Value2 = '8808605CS'
$ws.Range("D1309").Value2 = 'Torani Puremade Dark Chocolate Flavoring Sauce 64 fl. oz. - 4/Case'
$ws.Range("E1309").Value2 = ''' 8'
$ws.Range("H1309").Value2 = 67.98999999999999

# Row 1310
$ws.Range("A1310").Value2 = ' February 17, 2025'
$ws.Range("B1310").Value2 = '''107750625'
$ws.Range("C1310").Value2 = '110TWN05328'
$ws.Range("D1310").Value2 = 'Twinings Irish Breakfast Tea Bags - 20/Box'
$ws.Range("E1310").Value2 = ''' 12'
$ws.Range("H1310").Value2 = 4.39

# Row 1311
$ws.Range("A1311").Value2 = ' February 17, 2025'
$ws.Range("B1311").Value2 = '''107750625'
$ws.Range("C1311").Value2 = '711SPRNKLEPK'
$ws.Range("D1311").Value2 = 'Adourne Pink Sprinkles 10 lb.'
$ws.Range("E1311").Value2 = ''' 2'
$ws.Range("H1311").Value2 = 25.62

# Row 1312
$ws.Range("A1312").Value2 = ' February 17, 2025'
$ws.Range("B1312").Value2 = '''107750625'
$ws.Range("C1312").Value2 = '110TWN09183'
$ws.Range("D1312").Value2 = 'Twinings Earl Grey Tea Bags - 25/Box'
$ws.Range("E1312").Value2 = ''' 12'
$ws.Range("H1312").Value2 = 4.59

# Row 1313
$ws.Range("A1313").Value2 = ' February 17, 2025'
$ws.Range("B1313").Value2 = '''107750760'
$ws.Range("C1313").Value2 = '323CIRCL2018'
$ws.Range("D1313").Value2 = 'Lavex 2" Fluorescent Green Matte Paper Permanent Round Inventory Label - 500/Roll'
$ws.Range("E1313").Value2 = ''' 12'
$ws.Range("H1313").Value2 = 6.489999999999999

# Row 1314
$ws.Range("A1314").Value2 = ' February 18, 2025'
$ws.Range("B1314").Value2 = '''107786615'
$ws.Range("C1314").Value2 = '150BB6218N'
$ws.Range("D1314").Value2 = 'Durable Packaging BB6218N 18" x 6" x 3 1/4" Kraft Paper Windowed Bread Bag - 1000/Case'
$ws.Range("E1314").Value2 = ''' 2'
$ws.Range("H1314").Value2 = 104.99

# Row 1315
$ws.Range("A1315").Value2 = ' February 18, 2025'
$ws.Range("B1315").Value2 = '''107786615'
$ws.Range("C1315").Value2 = '43306LPIE300'
$ws.Range("D1315").Value2 = 'Choice 6" Clear Hinged Pie Container with Low Dome Lid - 300/Case'
$ws.Range("E1315").Value2 = ''' 4'
$ws.Range("H1315").Value2 = 49.99

# Row 1316
$ws.Range("A1316").Value2 = ' February 18, 2025'
$ws.Range("B1316").Value2 = '''107786615'
$ws.Range("C1316").Value2 = '''150300865'
$ws.Range("D1316").Value2 = 'Bagcraft Packaging 300865 EcoCraft 6" x 3 1/2" x 13 1/2" Dubl Panel® Artisan Bread Bag - 500/Case'
$ws.Range("E1316").Value2 = ''' 2'
$ws.Range("H1316").Value2 = 79.98999999999999

# Row 1317
$ws.Range("A1317").Value2 = ' February 18, 2025'
$ws.Range("B1317").Value2 = '''107786615'
$ws.Range("C1317").Value2 = '500CTOUT160'
$ws.Range("D1317").Value2 = 'Choice 160 oz. Beverage Take-Out Container - 20/Case'
$ws.Range("E1317").Value2 = ''' 3'
$ws.Range("H1317").Value2 = 94.99000000000001

# Row 1318
$ws.Range("A1318").Value2 = ' February 18, 2025'
$ws.Range("B1318").Value2 = '''107786615'
$ws.Range("C1318").Value2 = '433QLINERBL'
$ws.Range("D1318").Value2 = 'Baker''s Lane 16" x 24" Full Size Quilon® Coated Parchment Paper Bun / Sheet Pan Liner Sheet - 1000/Case'
$ws.Range("E1318").Value2 = ''' 6'
$ws.Range("H1318").Value2 = 45.99

# Row 1319
$ws.Range("A1319").Value2 = ' February 18, 2025'
$ws.Range("B1319").Value2 = '''107786615'
$ws.Range("C1319").Value2 = '433SLINERBL'
$ws.Range("D1319").Value2 = 'Baker''s Lane 16" x 24" Full Size Silicone Coated Parchment Paper Bun / Sheet Pan Liner Sheet - 1000/Case'
$ws.Range("E1319").Value2 = ''' 4'
$ws.Range("H1319").Value2 = 76.98999999999999

# Row 1320
$ws.Range("A1320").Value2 = ' February 19, 2025'
$ws.Range("B1320").Value2 = '''107836987'
$ws.Range("C1320").Value2 = '8808606CS'
$ws.Range("D1320").Value2 = 'Torani Puremade Caramel Flavoring Sauce 64 fl. oz. - 4/Case'
$ws.Range("E1320").Value2 = ''' 8'
$ws.Range("H1320").Value2 = 72.98999999999999

# Row 1321
$ws.Range("A1321").Value2 = ' February 23, 2025'
$ws.Range("B1321").Value2 = '''107965562'
$ws.Range("C1321").Value2 = '844CBBLKNBR35'
$ws.Range("D1321").Value2 = 'Lavex 3'' x 5'' Heavy-Duty Black Grease-Resistant Anti-Fatigue Closed-Cell Nitrile Rubber Floor Mat - 3/4" Thick'
$ws.Range("E1321").Value2 = ''' 2'
$ws.Range("H1321").Value2 = 92.98999999999999

# Row 1322
$ws.Range("A1322").Value2 = ' February 23, 2025'
$ws.Range("B1322").Value2 = '''107965562'
$ws.Range("C1322").Value2 = '697STW60'
$ws.Range("D1322").Value2 = 'Lavex 60" Wooden Mop Handle with Stirrup-Style End'
$ws.Range("E1322").Value2 = ''' 3'
$ws.Range("H1322").Value2 = 7.989999999999999

# Row 1323
$ws.Range("A1323").Value2 = ' February 23, 2025'
$ws.Range("B1323").Value2 = '''107965562'
$ws.Range("C1323").Value2 = '544SYPFR145K'
$ws.Range("D1323").Value2 = 'Monin Premium Toasted Marshmallow Flavoring Syrup 1 Liter - 4/Case'
$ws.Range("E1323").Value2 = ''' 2'
$ws.Range("H1323").Value2 = 37.99

# Row 1324
$ws.Range("A1324").Value2 = ' February 23, 2025'
$ws.Range("B1324").Value2 = '''107965562'
$ws.Range("C1324").Value2 = '544SYPFR023K'
$ws.Range("D1324").Value2 = 'Monin Premium Hazelnut Flavoring Syrup 1 Liter - 4/Case'
$ws.Range("E1324").Value2 = ''' 2'
$ws.Range("H1324").Value2 = 38.49

# Row 1325
$ws.Range("A1325").Value2 = ' February 23, 2025'
$ws.Range("B1325").Value2 = '''107965562'
$ws.Range("C1325").Value2 = '544SYPAR012KT'
$ws.Range("D1325").Value2 = 'Monin Premium Cinnamon Flavoring Syrup 750 mL - 12/Case'
$ws.Range("E1325").Value2 = ''' 1'
$ws.Range("H1325").Value2 = 86.98999999999999

# Row 1326
$ws.Range("A1326").Value2 = ' February 23, 2025'
$ws.Range("B1326").Value2 = '''107965562'
$ws.Range("C1326").Value2 = '544SYPAR258KT'
$ws.Range("D1326").Value2 = 'Monin Premium Butter Pecan Flavoring Syrup 750 mL - 12/Case'
$ws.Range("E1326").Value2 = ''' 2'
$ws.Range("H1326").Value2 = 86.98999999999999

# Row 1327
$ws.Range("A1327").Value2 = ' February 23, 2025'
$ws.Range("B1327").Value2 = '''107965562'
$ws.Range("C1327").Value2 = '110TWN08465KT'
$ws.Range("D1327").Value2 = 'Twinings Pomegranate & Raspberry Herbal Tea Bags - 120/Case'
$ws.Range("E1327").Value2 = ''' 2'
$ws.Range("H1327").Value2 = 21.99

# Row 1328
$ws.Range("A1328").Value2 = ' February 23, 2025'
$ws.Range("B1328").Value2 = '''107965562'
$ws.Range("C1328").Value2 = '182RRF8'
$ws.Range("D1328").Value2 = 'Choice 6 1/2" x 7 3/4" Plastic Food Bag on a Roll - 2000/Case'
$ws.Range("E1328").Value2 = ''' 6'
$ws.Range("H1328").Value2 = 19.29

# Row 1329
$ws.Range("A1329").Value2 = ' February 23, 2025'
$ws.Range("B1329").Value2 = '''107965562'
$ws.Range("C1329").Value2 = '110TWN05322KT'
$ws.Range("D1329").Value2 = 'Twinings Darjeeling Tea Bags - 120/Case'
$ws.Range("E1329").Value2 = ''' 2'
$ws.Range("H1329").Value2 = 21.99

# Row 1330
$ws.Range("A1330").Value2 = ' February 23, 2025'
$ws.Range("B1330").Value2 = '''107965562'
$ws.Range("C1330").Value2 = '8808607CS'
$ws.Range("D1330").Value2 = 'Torani Puremade Pumpkin Pie Flavoring Sauce 64 fl. oz. - 4/Case'
$ws.Range("E1330").Value2 = ''' 1'
$ws.Range("H1330").Value2 = 72.98999999999999

# Row 1331
$ws.Range("A1331").Value2 = ' February 26, 2025'
$ws.Range("B1331").Value2 = '''108078832'
$ws.Range("C1331").Value2 = '612H18A'
$ws.Range("D1331").Value2 = 'Durable Packaging High Dome Plastic Cover for 1/4 Sheet Cake Pan - 100/Case'
$ws.Range("E1331").Value2 = ''' 8'
$ws.Range("H1331").Value2 = 50.99

# Row 1332
$ws.Range("A1332").Value2 = ' February 26, 2025'
$ws.Range("B1332").Value2 = '''108078832'
$ws.Range("C1332").Value2 = '''612604245'
$ws.Range("D1332").Value2 = 'Durable Packaging 1/4 Sheet Foil Cake Pan - 100/Case'
$ws.Range("E1332").Value2 = ''' 8'
$ws.Range("H1332").Value2 = 48.99

# Row 1333
$ws.Range("A1333").Value2 = ' February 26, 2025'
$ws.Range("B1333").Value2 = '''108078832'
$ws.Range("C1333").Value2 = '111QUINTRORG'
$ws.Range("D1333").Value2 = 'Organic Tri-Color Quinoa - 25 lb.'
$ws.Range("E1333").Value2 = ''' 3'
$ws.Range("H1333").Value2 = 62.11000000000001

# Row 1334
$ws.Range("A1334").Value2 = ' February 26, 2025'
$ws.Range("B1334").Value2 = '''108082685'
$ws.Range("C1334").Value2 = '697JWMH60GN'
$ws.Range("D1334").Value2 = 'Lavex 60" Green Jaw Style Metal Mop Handle'
$ws.Range("E1334").Value2 = ''' 3'
$ws.Range("H1334").Value2 = 9.99

# Row 1335
$ws.Range("A1335").Value2 = ' February 26, 2025'
$ws.Range("B1335").Value2 = '''108082685'
$ws.Range("C1335").Value2 = '697QMH60GN'
$ws.Range("D1335").Value2 = 'Lavex 60" Green Quick Release Metal Mop Handle'
$ws.Range("E1335").Value2 = ''' 1'
$ws.Range("H1335").Value2 = 7.99

